$wb = $excel.ActiveWorkbook

# The task list lives on Sheet2 (the active sheet). Column F holds the
# "Status" for each task. Rows 9 and 10 ("Code Controller" / "Fix bug")
# were "In Process" / "Not started" - mark them both "Completed".
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F9").Value = "Completed"
$ws2.Range("F10").Value = "Completed"

# Move the active selection down to F11, matching where the cursor ended
# up after updating the two status cells.
$ws2.Range("F11").Select() | Out-Null
